$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet: insert a new row for 2022-Q3 at the top of the data table,
#    shift everything else down, and renumber the index column.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Insert a new blank row above the current row 2 (old 2022-Q1 row); this
# shifts every existing data row down by one and keeps their formatting.
$summary.Rows.Item(2).Insert()

# Row 2 lost its style during the blank insert (there was nothing above it
# in column A to inherit from) -- pull the formatting from row 3 which still
# carries the original "index" column styling.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
# Columns B-D in the data rows carry no explicit style (only column A does);
# the row-insert operation picked up the header's bold style for these
# cells, so reset them back to the plain/default formatting used by the
# other data rows.
$summary.Range("B2:D2").ClearFormats()

# Fill in the new 2022-Q3 summary row.
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 44
$summary.Cells.Item(2, 4).Value = 19.13

# Renumber the index column (A) for every row that was pushed down.
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(6, 1).Value = 4
$summary.Cells.Item(7, 1).Value = 5
$summary.Cells.Item(8, 1).Value = 6

# ---------------------------------------------------------------------------
# 2) Create the new "2022-Q3" holdings sheet. Clone the "2022-Q1" sheet (for
#    its layout/formatting) and place the clone right before it, then
#    overwrite its contents with the 2022-Q3 figures.
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Copy($q1)
$q3 = $wb.Worksheets.Item("2022-Q1 (2)")
$q3.Name = "2022-Q3"

# 2022-Q1 has 64 data rows (rows 2-65); 2022-Q3 only has 44 (rows 2-45).
# Drop the extra 20 rows so the sheet's used range matches.
$q3.Range("A46:H65").EntireRow.Delete()

# Columns B-G hold text in this workbook (fund code / name / size / position
# / weight / market value), even when the values look numeric -- force them
# to Text so values like "163417" or "19.13" are not re-interpreted as
# numbers.
$q3.Range("B2:G45").NumberFormat = "@"

$q3.Cells.Item(2, 1).Value = 0
$q3.Cells.Item(2, 2).Value = "163417"
$q3.Cells.Item(2, 3).Value = "兴全合宜灵活配置混合（LOF）A"
$q3.Cells.Item(2, 4).Value = "163.32"
$q3.Cells.Item(2, 5).Value = "92.65"
$q3.Cells.Item(2, 6).Value = "3.55"
$q3.Cells.Item(2, 7).Value = "5.7979"
$q3.Cells.Item(2, 8).Value = 5

$q3.Cells.Item(3, 1).Value = 1
$q3.Cells.Item(3, 2).Value = "002685"
$q3.Cells.Item(3, 3).Value = "中欧丰泓沪港深灵活配置混合A"
$q3.Cells.Item(3, 4).Value = "42.36"
$q3.Cells.Item(3, 5).Value = "92.77"
$q3.Cells.Item(3, 6).Value = "4.10"
$q3.Cells.Item(3, 7).Value = "1.7368"
$q3.Cells.Item(3, 8).Value = 9

$q3.Cells.Item(4, 1).Value = 2
$q3.Cells.Item(4, 2).Value = "166025"
$q3.Cells.Item(4, 3).Value = "中欧远见两年定期开放混合A"
$q3.Cells.Item(4, 4).Value = "44.75"
$q3.Cells.Item(4, 5).Value = "59.87"
$q3.Cells.Item(4, 6).Value = "2.94"
$q3.Cells.Item(4, 7).Value = "1.3156"
$q3.Cells.Item(4, 8).Value = 7

$q3.Cells.Item(5, 1).Value = 3
$q3.Cells.Item(5, 2).Value = "008378"
$q3.Cells.Item(5, 3).Value = "兴全社会价值三年持有期混合"
$q3.Cells.Item(5, 4).Value = "38.47"
$q3.Cells.Item(5, 5).Value = "91.22"
$q3.Cells.Item(5, 6).Value = "3.30"
$q3.Cells.Item(5, 7).Value = "1.2695"
$q3.Cells.Item(5, 8).Value = 6

$q3.Cells.Item(6, 1).Value = 4
$q3.Cells.Item(6, 2).Value = "012647"
$q3.Cells.Item(6, 3).Value = "中欧洞见一年持有混合"
$q3.Cells.Item(6, 4).Value = "33.23"
$q3.Cells.Item(6, 5).Value = "80.48"
$q3.Cells.Item(6, 6).Value = "2.60"
$q3.Cells.Item(6, 7).Value = "0.8640"
$q3.Cells.Item(6, 8).Value = 5

$q3.Cells.Item(7, 1).Value = 5
$q3.Cells.Item(7, 2).Value = "010723"
$q3.Cells.Item(7, 3).Value = "中欧价值成长混合A"
$q3.Cells.Item(7, 4).Value = "22.79"
$q3.Cells.Item(7, 5).Value = "86.71"
$q3.Cells.Item(7, 6).Value = "3.76"
$q3.Cells.Item(7, 7).Value = "0.8569"
$q3.Cells.Item(7, 8).Value = 9

$q3.Cells.Item(8, 1).Value = 6
$q3.Cells.Item(8, 2).Value = "009896"
$q3.Cells.Item(8, 3).Value = "广发港股通成长精选股票A"
$q3.Cells.Item(8, 4).Value = "18.30"
$q3.Cells.Item(8, 5).Value = "90.12"
$q3.Cells.Item(8, 6).Value = "4.23"
$q3.Cells.Item(8, 7).Value = "0.7741"
$q3.Cells.Item(8, 8).Value = 9

$q3.Cells.Item(9, 1).Value = 7
$q3.Cells.Item(9, 2).Value = "010947"
$q3.Cells.Item(9, 3).Value = "中欧嘉选混合A"
$q3.Cells.Item(9, 4).Value = "13.56"
$q3.Cells.Item(9, 5).Value = "85.30"
$q3.Cells.Item(9, 6).Value = "4.75"
$q3.Cells.Item(9, 7).Value = "0.6441"
$q3.Cells.Item(9, 8).Value = 5

$q3.Cells.Item(10, 1).Value = 8
$q3.Cells.Item(10, 2).Value = "011162"
$q3.Cells.Item(10, 3).Value = "博时港股通领先趋势混合A"
$q3.Cells.Item(10, 4).Value = "15.21"
$q3.Cells.Item(10, 5).Value = "90.30"
$q3.Cells.Item(10, 6).Value = "4.23"
$q3.Cells.Item(10, 7).Value = "0.6434"
$q3.Cells.Item(10, 8).Value = 10

$q3.Cells.Item(11, 1).Value = 9
$q3.Cells.Item(11, 2).Value = "002387"
$q3.Cells.Item(11, 3).Value = "工银沪港深股票A"
$q3.Cells.Item(11, 4).Value = "13.37"
$q3.Cells.Item(11, 5).Value = "83.30"
$q3.Cells.Item(11, 6).Value = "4.62"
$q3.Cells.Item(11, 7).Value = "0.6177"
$q3.Cells.Item(11, 8).Value = 4

$q3.Cells.Item(12, 1).Value = 10
$q3.Cells.Item(12, 2).Value = "010678"
$q3.Cells.Item(12, 3).Value = "中欧均衡成长混合A"
$q3.Cells.Item(12, 4).Value = "14.02"
$q3.Cells.Item(12, 5).Value = "87.25"
$q3.Cells.Item(12, 6).Value = "4.14"
$q3.Cells.Item(12, 7).Value = "0.5804"
$q3.Cells.Item(12, 8).Value = 8

$q3.Cells.Item(13, 1).Value = 11
$q3.Cells.Item(13, 2).Value = "005491"
$q3.Cells.Item(13, 3).Value = "兴全合宜灵活配置混合（LOF）C"
$q3.Cells.Item(13, 4).Value = "10.59"
$q3.Cells.Item(13, 5).Value = "92.65"
$q3.Cells.Item(13, 6).Value = "3.55"
$q3.Cells.Item(13, 7).Value = "0.3759"
$q3.Cells.Item(13, 8).Value = 5

$q3.Cells.Item(14, 1).Value = 12
$q3.Cells.Item(14, 2).Value = "513980"
$q3.Cells.Item(14, 3).Value = "景顺长城中证港股通科技ETF"
$q3.Cells.Item(14, 4).Value = "13.30"
$q3.Cells.Item(14, 5).Value = "98.24"
$q3.Cells.Item(14, 6).Value = "2.70"
$q3.Cells.Item(14, 7).Value = "0.3591"
$q3.Cells.Item(14, 8).Value = 10

$q3.Cells.Item(15, 1).Value = 13
$q3.Cells.Item(15, 2).Value = "159636"
$q3.Cells.Item(15, 3).Value = "工银瑞信国证港股通科技ETF"
$q3.Cells.Item(15, 4).Value = "8.42"
$q3.Cells.Item(15, 5).Value = "98.06"
$q3.Cells.Item(15, 6).Value = "3.72"
$q3.Cells.Item(15, 7).Value = "0.3132"
$q3.Cells.Item(15, 8).Value = 7

$q3.Cells.Item(16, 1).Value = 14
$q3.Cells.Item(16, 2).Value = "002686"
$q3.Cells.Item(16, 3).Value = "中欧丰泓沪港深灵活配置混合C"
$q3.Cells.Item(16, 4).Value = "7.40"
$q3.Cells.Item(16, 5).Value = "92.77"
$q3.Cells.Item(16, 6).Value = "4.10"
$q3.Cells.Item(16, 7).Value = "0.3034"
$q3.Cells.Item(16, 8).Value = 9

$q3.Cells.Item(17, 1).Value = 15
$q3.Cells.Item(17, 2).Value = "011856"
$q3.Cells.Item(17, 3).Value = "安信均衡成长18个月持有混合A"
$q3.Cells.Item(17, 4).Value = "5.25"
$q3.Cells.Item(17, 5).Value = "92.50"
$q3.Cells.Item(17, 6).Value = "5.64"
$q3.Cells.Item(17, 7).Value = "0.2961"
$q3.Cells.Item(17, 8).Value = 7

$q3.Cells.Item(18, 1).Value = 16
$q3.Cells.Item(18, 2).Value = "013991"
$q3.Cells.Item(18, 3).Value = "中欧港股通精选一年持有混合A"
$q3.Cells.Item(18, 4).Value = "6.69"
$q3.Cells.Item(18, 5).Value = "93.38"
$q3.Cells.Item(18, 6).Value = "4.00"
$q3.Cells.Item(18, 7).Value = "0.2676"
$q3.Cells.Item(18, 8).Value = 8

$q3.Cells.Item(19, 1).Value = 17
$q3.Cells.Item(19, 2).Value = "005241"
$q3.Cells.Item(19, 3).Value = "中欧时代智慧混合A"
$q3.Cells.Item(19, 4).Value = "6.95"
$q3.Cells.Item(19, 5).Value = "77.65"
$q3.Cells.Item(19, 6).Value = "3.69"
$q3.Cells.Item(19, 7).Value = "0.2565"
$q3.Cells.Item(19, 8).Value = 8

$q3.Cells.Item(20, 1).Value = 18
$q3.Cells.Item(20, 2).Value = "009897"
$q3.Cells.Item(20, 3).Value = "广发港股通成长精选股票C"
$q3.Cells.Item(20, 4).Value = "5.86"
$q3.Cells.Item(20, 5).Value = "90.12"
$q3.Cells.Item(20, 6).Value = "4.23"
$q3.Cells.Item(20, 7).Value = "0.2479"
$q3.Cells.Item(20, 8).Value = 9

$q3.Cells.Item(21, 1).Value = 19
$q3.Cells.Item(21, 2).Value = "013992"
$q3.Cells.Item(21, 3).Value = "中欧港股通精选一年持有混合C"
$q3.Cells.Item(21, 4).Value = "4.68"
$q3.Cells.Item(21, 5).Value = "93.38"
$q3.Cells.Item(21, 6).Value = "4.00"
$q3.Cells.Item(21, 7).Value = "0.1872"
$q3.Cells.Item(21, 8).Value = 8

$q3.Cells.Item(22, 1).Value = 20
$q3.Cells.Item(22, 2).Value = "011163"
$q3.Cells.Item(22, 3).Value = "博时港股通领先趋势混合C"
$q3.Cells.Item(22, 4).Value = "3.82"
$q3.Cells.Item(22, 5).Value = "90.30"
$q3.Cells.Item(22, 6).Value = "4.23"
$q3.Cells.Item(22, 7).Value = "0.1616"
$q3.Cells.Item(22, 8).Value = 10

$q3.Cells.Item(23, 1).Value = 21
$q3.Cells.Item(23, 2).Value = "011708"
$q3.Cells.Item(23, 3).Value = "中欧嘉益一年混合A"
$q3.Cells.Item(23, 4).Value = "4.44"
$q3.Cells.Item(23, 5).Value = "93.35"
$q3.Cells.Item(23, 6).Value = "3.49"
$q3.Cells.Item(23, 7).Value = "0.1550"
$q3.Cells.Item(23, 8).Value = 10

$q3.Cells.Item(24, 1).Value = 22
$q3.Cells.Item(24, 2).Value = "008891"
$q3.Cells.Item(24, 3).Value = "安信价值成长混合A"
$q3.Cells.Item(24, 4).Value = "2.69"
$q3.Cells.Item(24, 5).Value = "92.09"
$q3.Cells.Item(24, 6).Value = "5.64"
$q3.Cells.Item(24, 7).Value = "0.1517"
$q3.Cells.Item(24, 8).Value = 7

$q3.Cells.Item(25, 1).Value = 23
$q3.Cells.Item(25, 2).Value = "009880"
$q3.Cells.Item(25, 3).Value = "安信成长动力一年持有期混合"
$q3.Cells.Item(25, 4).Value = "1.91"
$q3.Cells.Item(25, 5).Value = "93.89"
$q3.Cells.Item(25, 6).Value = "6.39"
$q3.Cells.Item(25, 7).Value = "0.1220"
$q3.Cells.Item(25, 8).Value = 7

$q3.Cells.Item(26, 1).Value = 24
$q3.Cells.Item(26, 2).Value = "501021"
$q3.Cells.Item(26, 3).Value = "华宝标普香港上市中国中小盘指数（LOF）A"
$q3.Cells.Item(26, 4).Value = "4.19"
$q3.Cells.Item(26, 5).Value = "92.99"
$q3.Cells.Item(26, 6).Value = "2.28"
$q3.Cells.Item(26, 7).Value = "0.0955"
$q3.Cells.Item(26, 8).Value = 3

$q3.Cells.Item(27, 1).Value = 25
$q3.Cells.Item(27, 2).Value = "008892"
$q3.Cells.Item(27, 3).Value = "安信价值成长混合C"
$q3.Cells.Item(27, 4).Value = "1.53"
$q3.Cells.Item(27, 5).Value = "92.09"
$q3.Cells.Item(27, 6).Value = "5.64"
$q3.Cells.Item(27, 7).Value = "0.0863"
$q3.Cells.Item(27, 8).Value = 7

$q3.Cells.Item(28, 1).Value = 26
$q3.Cells.Item(28, 2).Value = "005242"
$q3.Cells.Item(28, 3).Value = "中欧时代智慧混合C"
$q3.Cells.Item(28, 4).Value = "2.29"
$q3.Cells.Item(28, 5).Value = "77.65"
$q3.Cells.Item(28, 6).Value = "3.69"
$q3.Cells.Item(28, 7).Value = "0.0845"
$q3.Cells.Item(28, 8).Value = 8

$q3.Cells.Item(29, 1).Value = 27
$q3.Cells.Item(29, 2).Value = "007101"
$q3.Cells.Item(29, 3).Value = "中欧远见两年定期开放混合C"
$q3.Cells.Item(29, 4).Value = "2.79"
$q3.Cells.Item(29, 5).Value = "59.87"
$q3.Cells.Item(29, 6).Value = "2.94"
$q3.Cells.Item(29, 7).Value = "0.0820"
$q3.Cells.Item(29, 8).Value = 7

$q3.Cells.Item(30, 1).Value = 28
$q3.Cells.Item(30, 2).Value = "012379"
$q3.Cells.Item(30, 3).Value = "创金合信港股互联网3个月持有期混合（QDII）A"
$q3.Cells.Item(30, 4).Value = "2.81"
$q3.Cells.Item(30, 5).Value = "87.48"
$q3.Cells.Item(30, 6).Value = "2.42"
$q3.Cells.Item(30, 7).Value = "0.0680"
$q3.Cells.Item(30, 8).Value = 10

$q3.Cells.Item(31, 1).Value = 29
$q3.Cells.Item(31, 2).Value = "011709"
$q3.Cells.Item(31, 3).Value = "中欧嘉益一年混合C"
$q3.Cells.Item(31, 4).Value = "1.87"
$q3.Cells.Item(31, 5).Value = "93.35"
$q3.Cells.Item(31, 6).Value = "3.49"
$q3.Cells.Item(31, 7).Value = "0.0653"
$q3.Cells.Item(31, 8).Value = 10

$q3.Cells.Item(32, 1).Value = 30
$q3.Cells.Item(32, 2).Value = "010724"
$q3.Cells.Item(32, 3).Value = "中欧价值成长混合C"
$q3.Cells.Item(32, 4).Value = "1.57"
$q3.Cells.Item(32, 5).Value = "86.71"
$q3.Cells.Item(32, 6).Value = "3.76"
$q3.Cells.Item(32, 7).Value = "0.0590"
$q3.Cells.Item(32, 8).Value = 9

$q3.Cells.Item(33, 1).Value = 31
$q3.Cells.Item(33, 2).Value = "011924"
$q3.Cells.Item(33, 3).Value = "嘉实港股互联网产业核心资产混合A"
$q3.Cells.Item(33, 4).Value = "1.20"
$q3.Cells.Item(33, 5).Value = "87.88"
$q3.Cells.Item(33, 6).Value = "4.06"
$q3.Cells.Item(33, 7).Value = "0.0487"
$q3.Cells.Item(33, 8).Value = 8

$q3.Cells.Item(34, 1).Value = 32
$q3.Cells.Item(34, 2).Value = "007512"
$q3.Cells.Item(34, 3).Value = "工银沪港深股票C"
$q3.Cells.Item(34, 4).Value = "0.86"
$q3.Cells.Item(34, 5).Value = "83.30"
$q3.Cells.Item(34, 6).Value = "4.62"
$q3.Cells.Item(34, 7).Value = "0.0397"
$q3.Cells.Item(34, 8).Value = 4

$q3.Cells.Item(35, 1).Value = 33
$q3.Cells.Item(35, 2).Value = "010948"
$q3.Cells.Item(35, 3).Value = "中欧嘉选混合C"
$q3.Cells.Item(35, 4).Value = "0.78"
$q3.Cells.Item(35, 5).Value = "85.30"
$q3.Cells.Item(35, 6).Value = "4.75"
$q3.Cells.Item(35, 7).Value = "0.0370"
$q3.Cells.Item(35, 8).Value = 5

$q3.Cells.Item(36, 1).Value = 34
$q3.Cells.Item(36, 2).Value = "513160"
$q3.Cells.Item(36, 3).Value = "银华恒生港股通中国科技ETF"
$q3.Cells.Item(36, 4).Value = "0.43"
$q3.Cells.Item(36, 5).Value = "92.45"
$q3.Cells.Item(36, 6).Value = "7.82"
$q3.Cells.Item(36, 7).Value = "0.0336"
$q3.Cells.Item(36, 8).Value = 6

$q3.Cells.Item(37, 1).Value = 35
$q3.Cells.Item(37, 2).Value = "010679"
$q3.Cells.Item(37, 3).Value = "中欧均衡成长混合C"
$q3.Cells.Item(37, 4).Value = "0.74"
$q3.Cells.Item(37, 5).Value = "87.25"
$q3.Cells.Item(37, 6).Value = "4.14"
$q3.Cells.Item(37, 7).Value = "0.0306"
$q3.Cells.Item(37, 8).Value = 8

$q3.Cells.Item(38, 1).Value = 36
$q3.Cells.Item(38, 2).Value = "012380"
$q3.Cells.Item(38, 3).Value = "创金合信港股互联网3个月持有期混合（QDII）C"
$q3.Cells.Item(38, 4).Value = "0.96"
$q3.Cells.Item(38, 5).Value = "87.48"
$q3.Cells.Item(38, 6).Value = "2.42"
$q3.Cells.Item(38, 7).Value = "0.0232"
$q3.Cells.Item(38, 8).Value = 10

$q3.Cells.Item(39, 1).Value = 37
$q3.Cells.Item(39, 2).Value = "003413"
$q3.Cells.Item(39, 3).Value = "华泰柏瑞新经济沪港深混合"
$q3.Cells.Item(39, 4).Value = "0.42"
$q3.Cells.Item(39, 5).Value = "86.45"
$q3.Cells.Item(39, 6).Value = "5.39"
$q3.Cells.Item(39, 7).Value = "0.0226"
$q3.Cells.Item(39, 8).Value = 8

$q3.Cells.Item(40, 1).Value = 38
$q3.Cells.Item(40, 2).Value = "011925"
$q3.Cells.Item(40, 3).Value = "嘉实港股互联网产业核心资产混合C"
$q3.Cells.Item(40, 4).Value = "0.41"
$q3.Cells.Item(40, 5).Value = "87.88"
$q3.Cells.Item(40, 6).Value = "4.06"
$q3.Cells.Item(40, 7).Value = "0.0166"
$q3.Cells.Item(40, 8).Value = 8

$q3.Cells.Item(41, 1).Value = 39
$q3.Cells.Item(41, 2).Value = "011857"
$q3.Cells.Item(41, 3).Value = "安信均衡成长18个月持有混合C"
$q3.Cells.Item(41, 4).Value = "0.29"
$q3.Cells.Item(41, 5).Value = "92.50"
$q3.Cells.Item(41, 6).Value = "5.64"
$q3.Cells.Item(41, 7).Value = "0.0164"
$q3.Cells.Item(41, 8).Value = 7

$q3.Cells.Item(42, 1).Value = 40
$q3.Cells.Item(42, 2).Value = "012990"
$q3.Cells.Item(42, 3).Value = "天弘国证港股通50指数C"
$q3.Cells.Item(42, 4).Value = "0.27"
$q3.Cells.Item(42, 5).Value = "93.05"
$q3.Cells.Item(42, 6).Value = "2.73"
$q3.Cells.Item(42, 7).Value = "0.0074"
$q3.Cells.Item(42, 8).Value = 10

$q3.Cells.Item(43, 1).Value = 41
$q3.Cells.Item(43, 2).Value = "004321"
$q3.Cells.Item(43, 3).Value = "前海开源沪港深强国产业灵活配置混合"
$q3.Cells.Item(43, 4).Value = "0.11"
$q3.Cells.Item(43, 5).Value = "78.52"
$q3.Cells.Item(43, 6).Value = "5.41"
$q3.Cells.Item(43, 7).Value = "0.0060"
$q3.Cells.Item(43, 8).Value = 7

$q3.Cells.Item(44, 1).Value = 42
$q3.Cells.Item(44, 2).Value = "006127"
$q3.Cells.Item(44, 3).Value = "华宝标普香港上市中国中小盘指数（LOF）C"
$q3.Cells.Item(44, 4).Value = "0.24"
$q3.Cells.Item(44, 5).Value = "92.99"
$q3.Cells.Item(44, 6).Value = "2.28"
$q3.Cells.Item(44, 7).Value = "0.0055"
$q3.Cells.Item(44, 8).Value = 3

$q3.Cells.Item(45, 1).Value = 43
$q3.Cells.Item(45, 2).Value = "012989"
$q3.Cells.Item(45, 3).Value = "天弘国证港股通50指数A"
$q3.Cells.Item(45, 4).Value = "0.02"
$q3.Cells.Item(45, 5).Value = "93.05"
$q3.Cells.Item(45, 6).Value = "2.73"
$q3.Cells.Item(45, 7).Value = "0.0005"
$q3.Cells.Item(45, 8).Value = 10

Write-Host "2022-Q3 sheet populated"
